# Delete the empty Sheet2 and Sheet3 worksheets, keeping only Sheet1,
# and rename Sheet1 to "Product_Components".

$wb = $excel.ActiveWorkbook

# Turn off alerts so deleting a sheet doesn't prompt a confirmation dialog.
$excel.DisplayAlerts = $false

$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Delete()

$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Delete()

$excel.DisplayAlerts = $true

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "Product_Components"
